# "added 4wk low sales check"
# Update the per-week forecast figures (MyForecast / Inventory Coverage /
# Seasonality Index) on the "Forecast Comparison" sheet, and the derived
# roll-up figures on the "Summary" sheet, to reflect the new forecast run
# that includes the 4-week low sales check.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------
# Columns: D = MyForecast, H = Inventory Coverage, L = Seasonality Index

# Row 2 (W10)
$wsForecast.Range("D2").Value = 36
$wsForecast.Range("H2").Value = 7.66
$wsForecast.Range("L2").Value = 1.11

# Row 3 (W11)
$wsForecast.Range("D3").Value = 36
$wsForecast.Range("H3").Value = 6.62
$wsForecast.Range("L3").Value = 1.1

# Row 4 (W12)
$wsForecast.Range("D4").Value = 36
$wsForecast.Range("H4").Value = 5.59
$wsForecast.Range("L4").Value = 1.19

# Row 5 (W13) - MyForecast unchanged
$wsForecast.Range("H5").Value = 4.57
$wsForecast.Range("L5").Value = 1.16

# Row 6 (W14)
$wsForecast.Range("D6").Value = 37
$wsForecast.Range("H6").Value = 3.55
$wsForecast.Range("L6").Value = 0.92

# Row 7 (W15)
$wsForecast.Range("D7").Value = 37
$wsForecast.Range("H7").Value = 2.52
$wsForecast.Range("L7").Value = 0.86

# Row 8 (W16)
$wsForecast.Range("D8").Value = 37
$wsForecast.Range("H8").Value = 1.52
$wsForecast.Range("L8").Value = 1.06

# Row 9 (W17) - MyForecast unchanged
$wsForecast.Range("H9").Value = 0.52
$wsForecast.Range("L9").Value = 0.95

# Row 10 (W18)
$wsForecast.Range("D10").Value = 37
$wsForecast.Range("L10").Value = 0.84

# Row 11 (W19)
$wsForecast.Range("D11").Value = 38
$wsForecast.Range("L11").Value = 0.92

# Row 12 (W20)
$wsForecast.Range("D12").Value = 38
$wsForecast.Range("L12").Value = 0.91

# Row 13 (W21)
$wsForecast.Range("D13").Value = 38
$wsForecast.Range("L13").Value = 1.18

# Row 14 (W22)
$wsForecast.Range("D14").Value = 38
$wsForecast.Range("L14").Value = 0.95

# Row 15 (W23)
$wsForecast.Range("D15").Value = 39
$wsForecast.Range("L15").Value = 0.8100000000000001

# Row 16 (W24)
$wsForecast.Range("D16").Value = 39
$wsForecast.Range("L16").Value = 1.18

# Row 17 (W25)
$wsForecast.Range("D17").Value = 39
$wsForecast.Range("L17").Value = 0.83

# --- Summary sheet ---------------------------------------------------------
# These "numeric-looking" totals are stored as text in the workbook, so the
# cells are pre-formatted as Text before assignment to stop Excel from
# auto-coercing the typed digits into a numeric value.
foreach ($addr in @("B9", "B10", "B11", "B12", "B14")) {
    $wsSummary.Range($addr).NumberFormat = "@"
}
$wsSummary.Range("B9").Value  = "606"
$wsSummary.Range("B10").Value = "296"
$wsSummary.Range("B11").Value = "146"
$wsSummary.Range("B12").Value = "40"
$wsSummary.Range("B14").Value = "36"
